# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# Column G (header "K") holds the recalculated strike-count values
# (s_vals). Recompute and write the new values for rows 2-64.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @(3,0,0,0,0,2,1,0,2,1,2,1,0,3,1,0,4,3,0,0,1,0,0,1,1,1,0,1,2,2,3,1,0,1,1,0,2,0,0,4,0,0,0,0,2,1,1,1,1,1,1,1,1,2,0,0,0,1,3,1,0,2,1)

$startRow = 2
for ($i = 0; $i -lt $newK.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 7).Value = $newK[$i]
}
